$wb = $excel.ActiveWorkbook
$sheet0 = $wb.Worksheets.Item("Sheet0")

# Remove the A2 value from Sheet0 (leaves just A1)
$sheet0.Range("A2").ClearContents()
$sheet0.PageSetup.Orientation = 1
$sheet0.Columns.Item(1).ColumnWidth = 15.166666666666666

# Insert a new worksheet right after "Sheet0"
$newSheet = $wb.Worksheets.Add($null, $sheet0)
$newSheet.Name = "Txn of Regression UAT2"

# Write cells in the precise order that reproduces the original shared-string allocation order.
$newSheet.Range("C2").Value = "AFSPKCRACC"
$newSheet.Range("C3").Value = "Current Account"
$newSheet.Range("C4").Value = " CA AKK"
$newSheet.Range("C5").Value = " CD BBA"

$newSheet.Range("B1").Value = "Category"
$newSheet.Range("C1").Value = "Account type Name"
$newSheet.Range("D1").Value = "Customer ID"

$newSheet.Range("C8").Value = "AlfalahKifayat"

$newSheet.Range("E1").Value = "Customer Type"

$newSheet.Range("E2").Value = "Minor"
$newSheet.Range("E3").Value = "Staff"
$newSheet.Range("E4").Value = "Individuals- Business / Self-Employ"
$newSheet.Range("E5").Value = "Individuals - Salaried"

$newSheet.Range("E8").Value = " Individuals - Others"

# Remaining duplicate text values (reuse existing shared strings, order doesn't add new entries)
$newSheet.Range("C6").Value = "AFSPKCRACC"
$newSheet.Range("C7").Value = "Current Account"
$newSheet.Range("E6").Value = "Minor"
$newSheet.Range("E7").Value = "Staff"

$newSheet.Range("A1").Value = "Account Number"

# Numeric columns
$newSheet.Range("A2").Value = 1008784258
$newSheet.Range("A3").Value = 1008784259
$newSheet.Range("A4").Value = 1008784262
$newSheet.Range("A5").Value = 1008784263
$newSheet.Range("A6").Value = 1008784270
$newSheet.Range("A7").Value = 1008784271
$newSheet.Range("A8").Value = 1008784273

$newSheet.Range("B2").Value = 1070
$newSheet.Range("B3").Value = 1001
$newSheet.Range("B4").Value = 1150
$newSheet.Range("B5").Value = 1005
$newSheet.Range("B6").Value = 1070
$newSheet.Range("B7").Value = 1001
$newSheet.Range("B8").Value = 6005

$newSheet.Range("D2").Value = 17866752
$newSheet.Range("D3").Value = 17866753
$newSheet.Range("D4").Value = 17866756
$newSheet.Range("D5").Value = 17866759
$newSheet.Range("D6").Value = 17866761
$newSheet.Range("D7").Value = 17866762
$newSheet.Range("D8").Value = 17866764

# Column widths (best-fit approximations of the original autofit widths)
$newSheet.Columns.Item(1).ColumnWidth = 15.166666666666666
$newSheet.Columns.Item(3).ColumnWidth = 17.666666666666668
$newSheet.Columns.Item(4).ColumnWidth = 11
$newSheet.Columns.Item(5).ColumnWidth = 31.833333333333332

[void]$newSheet.Range("J10").Select()

# Keep Sheet0 as the active sheet/tab
$sheet0.Activate()
[void]$sheet0.Range("A2:B8").Select()
